$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Valor Mora" total
$ws.Range("E11").Value = 113880

# Update "Cant. Trabajadores" count
$ws.Range("C13").Value = 1

# Update totals for the existing worker row (row 16)
$ws.Range("G16").Value = 1423500

# Row 17 previously held a second worker (DEIVIS ESPINOSA BALLESTAS, period 2506).
# It is replaced with another period entry (2508) for the same worker as row 16
# (JULIO CESAR TEHERAN SANTAMARIA), consistent with "Cant. Trabajadores" = 1.
$ws.Range("C17").Value = "9296727"
$ws.Range("D17").Value = "JULIO CESAR TEHERAN SANTAMARIA"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940

$wb.Save()
